$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = "'57210822856"
$ws.Range("B9").Value = 'T00000243'
$ws.Range("C9").Value = 'SONIA HELENA CONTRERAS ORTIZ'
$ws.Range("D9").Value = 'ESCUELA DE INGENIERÍA, ARQUITECTURA Y DISEÑO'
$ws.Range("E9").Value = 'Contreras-Ortiz, Sonia H.'
$ws.Range("F9").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57210822856'
$ws.Range("G9").Value = 'Contreras-Ortiz, S.H. (57210822856)'
$ws.Range("I9").Value = 5
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("S9").Value = 8
$ws.Range("U9").Value = 0
$ws.Range("X9").Value = 6
$ws.Range("Y9").Value = 2
$ws.Range("Z9").Value = 0
$ws.Range("AA9").Value = 0
$ws.Range("AC9").Value = 2
$ws.Range("AD9").Value = 3
$ws.Range("AE9").Value = 4
$ws.Range("AF9").Value = 0
$ws.Range("AH9").Value = 1
$ws.Range("AI9").Value = 7
$ws.Range("AJ9").Value = 3

# Row 10
$ws.Range("A10").Value = "'56674579200"
$ws.Range("B10").Value = 'T00020729'
$ws.Range("C10").Value = 'ROSA LEONOR ACEVEDO BARRIOS'
$ws.Range("D10").Value = 'CIENCIAS BÁSICAS'
$ws.Range("E10").Value = 'Acevedo-Barrios, R. L.'
$ws.Range("F10").Value = 'https://www.scopus.com/authid/detail.uri?authorId=56674579200'
$ws.Range("G10").Value = 'Acevedo-Barrios, R. (56674579200)'
$ws.Range("I10").Value = 7
$ws.Range("J10").Value = 2
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 2
$ws.Range("P10").Value = 1
$ws.Range("S10").Value = 5
$ws.Range("U10").Value = 1
$ws.Range("X10").Value = 3
$ws.Range("Y10").Value = 0
$ws.Range("Z10").Value = 2
$ws.Range("AA10").Value = 1
$ws.Range("AC10").Value = 1
$ws.Range("AD10").Value = 2
$ws.Range("AE10").Value = 2
$ws.Range("AF10").Value = 3
$ws.Range("AH10").Value = 2
$ws.Range("AI10").Value = 5
$ws.Range("AJ10").Value = 4

# Row 18
$ws.Range("A18").Value = "'57196040759"
$ws.Range("B18").Value = 'T00052721'
$ws.Range("C18").Value = 'PEDRO VAZQUEZ MIRAZ'
$ws.Range("D18").Value = 'ESCUELA DE NEGOCIOS, LEYES Y SOCIEDAD'
$ws.Range("E18").Value = 'Vázquez-Miraz, Pedro'
$ws.Range("F18").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57196040759'
$ws.Range("G18").Value = 'Vázquez-Miraz, Pedro (57196040759)'
$ws.Range("I18").Value = 3
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("O18").Value = 1
$ws.Range("Q18").Value = 3
$ws.Range("R18").Value = 1
$ws.Range("S18").Value = 1
$ws.Range("V18").Value = 0
$ws.Range("X18").Value = 0
$ws.Range("Z18").Value = 1
$ws.Range("AA18").Value = 2
$ws.Range("AC18").Value = 0
$ws.Range("AD18").Value = 1
$ws.Range("AE18").Value = 1
$ws.Range("AF18").Value = 5
$ws.Range("AG18").Value = 6
$ws.Range("AH18").Value = 2
$ws.Range("AI18").Value = 2
$ws.Range("AJ18").Value = 4
$ws.Range("AK18").Value = 4

# Row 19
$ws.Range("A19").Value = "'56581610900"
$ws.Range("B19").Value = 'T00007524'
$ws.Range("C19").Value = 'JUAN GABRIEL FAJARDO CUADRO'
$ws.Range("D19").Value = 'ESCUELA DE INGENIERÍA, ARQUITECTURA Y DISEÑO'
$ws.Range("E19").Value = 'Fajardo, Juan Gabriel'
$ws.Range("F19").Value = 'https://www.scopus.com/authid/detail.uri?authorId=56581610900'
$ws.Range("G19").Value = 'Fajardo-Cuadro, Juan Gabriel (56581610900)'
$ws.Range("Q19").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("S19").Value = 2
$ws.Range("U19").Value = 0
$ws.Range("V19").Value = 1
$ws.Range("W19").Value = 0
$ws.Range("X19").Value = 3
$ws.Range("Y19").Value = 0
$ws.Range("AA19").Value = 1
$ws.Range("AC19").Value = 2
$ws.Range("AF19").Value = 1
$ws.Range("AG19").Value = 3
$ws.Range("AH19").Value = 1
$ws.Range("AI19").Value = 4
$ws.Range("AJ19").Value = 2

# Row 20
$ws.Range("A20").Value = "'57219403758"
$ws.Range("B20").Value = 'T00028098'
$ws.Range("C20").Value = 'MARY JUDITH ARIAS TAPIA'
$ws.Range("D20").Value = 'ESCUELA DE INGENIERÍA, ARQUITECTURA Y DISEÑO'
$ws.Range("E20").Value = 'Arias Tapia, Mary Judith'
$ws.Range("F20").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57219403758'
$ws.Range("G20").Value = 'Arias-Tapia, Mary Judith (57219403758)'
$ws.Range("H20").Value = 3
$ws.Range("J20").Value = 1
$ws.Range("N20").Value = 0
$ws.Range("P20").Value = 1
$ws.Range("S20").Value = 0
$ws.Range("X20").Value = 0
$ws.Range("Y20").Value = 2
$ws.Range("AA20").Value = 1
$ws.Range("AC20").Value = 0
$ws.Range("AD20").Value = 1
$ws.Range("AF20").Value = 1
$ws.Range("AG20").Value = 1
$ws.Range("AI20").Value = 0
$ws.Range("AJ20").Value = 3
$ws.Range("AK20").Value = 3

# Row 21
$ws.Range("A21").Value = "'57220077867"
$ws.Range("B21").Value = 'T00000060'
$ws.Range("C21").Value = 'VILMA VIVIANA OJEDA CAICEDO'
$ws.Range("D21").Value = 'CIENCIAS BÁSICAS'
$ws.Range("E21").Value = 'Ojeda-Caicedo, Vilma Viviana'
$ws.Range("F21").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57220077867'
$ws.Range("G21").Value = 'Ojeda-Caicedo, Vilma V. (57220077867)'
$ws.Range("I21").Value = 3
$ws.Range("N21").Value = 1
$ws.Range("P21").Value = 0
$ws.Range("S21").Value = 5
$ws.Range("U21").Value = 0
$ws.Range("X21").Value = 3
$ws.Range("Y21").Value = 0
$ws.Range("AB21").Value = 0
$ws.Range("AC21").Value = 1
$ws.Range("AE21").Value = 1
$ws.Range("AF21").Value = 0
$ws.Range("AG21").Value = 5
$ws.Range("AH21").Value = 0
$ws.Range("AI21").Value = 3
$ws.Range("AJ21").Value = 0
$ws.Range("AK21").Value = 1

# Row 22
$ws.Range("A22").Value = "'57200615582"
$ws.Range("B22").Value = 'T00040553'
$ws.Range("C22").Value = 'MILTON CESAR GUERRERO PAJARO'
$ws.Range("D22").Value = 'ESCUELA DE INGENIERÍA, ARQUITECTURA Y DISEÑO'
$ws.Range("E22").Value = 'Guerrero, Milton'
$ws.Range("F22").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57200615582'
$ws.Range("G22").Value = 'Guerrero, Milton (57200615582)'
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("Q22").Value = 0
$ws.Range("U22").Value = 1
$ws.Range("W22").Value = 0
$ws.Range("AA22").Value = 0
$ws.Range("AB22").Value = 1
$ws.Range("AD22").Value = 0
$ws.Range("AE22").Value = 0
$ws.Range("AF22").Value = 1
$ws.Range("AG22").Value = 1

# Row 23
$ws.Range("A23").Value = "'57206773929"
$ws.Range("B23").Value = 'T00055756'
$ws.Range("C23").Value = 'TANIA LUCIA COBOS COBOS'
$ws.Range("D23").Value = 'ESCUELA DE TRANSFORMACIÓN DIGITAL'
$ws.Range("E23").Value = 'Cobos, Tania Lucía'
$ws.Range("F23").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57206773929'
$ws.Range("G23").Value = 'Cobos, Tania Lucía (57206773929)'
$ws.Range("I23").Value = 3
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("U23").Value = 1
$ws.Range("W23").Value = 0
$ws.Range("AA23").Value = 0
$ws.Range("AB23").Value = 0
$ws.Range("AC23").Value = 0
$ws.Range("AD23").Value = 0
$ws.Range("AE23").Value = 0
$ws.Range("AF23").Value = 1
$ws.Range("AG23").Value = 1
$ws.Range("AJ23").Value = 0
$ws.Range("AK23").Value = 0

# Row 24
$ws.Range("A24").Value = "'57758796500"
$ws.Range("B24").Value = 'T00000058'
$ws.Range("C24").Value = 'JORGE LUIS MUÑIZ OLITE'
$ws.Range("E24").Value = 'Muñiz-Olite, Jorge Luis'
$ws.Range("F24").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57758796500'
$ws.Range("G24").Value = 'Olite, Jorge Luis Muñiz (57758796500)'
$ws.Range("I24").Value = 4
$ws.Range("S24").Value = 0
$ws.Range("U24").Value = 1
$ws.Range("V24").Value = 1
$ws.Range("X24").Value = 0
$ws.Range("Z24").Value = 1
$ws.Range("AA24").Value = 1
$ws.Range("AC24").Value = 0
$ws.Range("AG24").Value = 1
$ws.Range("AH24").Value = 2
$ws.Range("AI24").Value = 2
$ws.Range("AJ24").Value = 2

# Row 25
$ws.Range("A25").Value = "'57192930752"
$ws.Range("B25").Value = 'T00019649'
$ws.Range("C25").Value = 'KAROL PATRICIA GUTIERREZ RUIZ'
$ws.Range("E25").Value = 'Gutiérrez-Ruiz, Karol Patricia'
$ws.Range("F25").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57192930752'
$ws.Range("G25").Value = 'Gutiérrez-Ruiz, K. (57192930752)'
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 1
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 1
$ws.Range("O25").Value = 1
$ws.Range("R25").Value = 1
$ws.Range("W25").Value = 1
$ws.Range("Y25").Value = 0
$ws.Range("AA25").Value = 2
$ws.Range("AB25").Value = 2
$ws.Range("AC25").Value = 1
$ws.Range("AD25").Value = 3
$ws.Range("AE25").Value = 3
$ws.Range("AF25").Value = 2
$ws.Range("AG25").Value = 2
$ws.Range("AH25").Value = 1
$ws.Range("AI25").Value = 1
$ws.Range("AJ25").Value = 4
$ws.Range("AK25").Value = 5

# Row 26
$ws.Range("A26").Value = "'57350116000"
$ws.Range("B26").Value = 'T00055760'
$ws.Range("C26").Value = 'JEOVANNY DE JESUS MUENTES ACEVEDO'
$ws.Range("D26").Value = 'CIENCIAS BÁSICAS'
$ws.Range("E26").Value = 'Acevedo, Jeovanny Muentes'
$ws.Range("F26").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57222990626'
$ws.Range("G26").Value = 'Muentes, J. (57350116000)'
$ws.Range("I26").Value = 2
$ws.Range("J26").Value = 1
$ws.Range("L26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 1
$ws.Range("U26").Value = 0
$ws.Range("W26").Value = 1
$ws.Range("Y26").Value = 1
$ws.Range("AD26").Value = 2
$ws.Range("AE26").Value = 2
$ws.Range("AF26").Value = 2
$ws.Range("AG26").Value = 2
$ws.Range("AJ26").Value = 2
$ws.Range("AK26").Value = 2

# Row 27
$ws.Range("A27").Value = "'56801043600"
$ws.Range("B27").Value = 'T00069460'
$ws.Range("C27").Value = 'ARGEMIRO PALENCIA DIAZ'
$ws.Range("D27").Value = 'ESCUELA DE INGENIERÍA, ARQUITECTURA Y DISEÑO'
$ws.Range("E27").Value = 'Palencia Díaz, Argemiro'
$ws.Range("F27").Value = 'https://www.scopus.com/authid/detail.uri?authorId=56801043600'
$ws.Range("G27").Value = 'Día, A. Palencia (56801043600)'
$ws.Range("I27").Value = 3
$ws.Range("O27").Value = 3
$ws.Range("P27").Value = 2
$ws.Range("V27").Value = 0
$ws.Range("Z27").Value = 0
$ws.Range("AF27").Value = 6
$ws.Range("AG27").Value = 6
$ws.Range("AH27").Value = 1
$ws.Range("AI27").Value = 1
$ws.Range("AJ27").Value = 1
$ws.Range("AK27").Value = 1

# Row 28
$ws.Range("A28").Value = "'56380539800"
$ws.Range("B28").Value = 'T00009384'
$ws.Range("C28").Value = 'DANIEL TORO GONZALEZ'
$ws.Range("D28").Value = 'ESCUELA DE NEGOCIOS, LEYES Y SOCIEDAD'
$ws.Range("E28").Value = 'Toro-González, Daniel'
$ws.Range("F28").Value = 'https://www.scopus.com/authid/detail.uri?authorId=56380539800'
$ws.Range("G28").Value = 'Toro Gonzalez, Daniel (56380539800)'
$ws.Range("I28").Value = 1
$ws.Range("R28").Value = 0
$ws.Range("U28").Value = 0
$ws.Range("Y28").Value = 1
$ws.Range("AF28").Value = 0
$ws.Range("AG28").Value = 0
$ws.Range("AH28").Value = 0
$ws.Range("AI28").Value = 0
$ws.Range("AJ28").Value = 1
$ws.Range("AK28").Value = 1

# Row 29
$ws.Range("A29").Value = "'57188841051"
$ws.Range("B29").Value = 'T00053267'
$ws.Range("C29").Value = 'EDISSON CHAVARRO MESA'
$ws.Range("D29").Value = 'CIENCIAS BÁSICAS'
$ws.Range("E29").Value = 'Chavarro-Mesa, Edisson'
$ws.Range("F29").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57188841051'
$ws.Range("G29").Value = 'Chavarro-Mesa, Edisson (57188841051)'
$ws.Range("I29").Value = 2
$ws.Range("L29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("Q29").Value = 1
$ws.Range("S29").Value = 4
$ws.Range("V29").Value = 0
$ws.Range("X29").Value = 5
$ws.Range("Z29").Value = 0
$ws.Range("AA29").Value = 0
$ws.Range("AD29").Value = 0
$ws.Range("AE29").Value = 0
$ws.Range("AG29").Value = 5
$ws.Range("AH29").Value = 0
$ws.Range("AI29").Value = 5
$ws.Range("AJ29").Value = 0
$ws.Range("AK29").Value = 2

# Row 30
$ws.Range("A30").Value = "'57202285682"
$ws.Range("B30").Value = 'T00008436'
$ws.Range("C30").Value = 'EDWIN ALEXANDER PUERTAS DEL CASTILLO'
$ws.Range("D30").Value = 'ESCUELA DE TRANSFORMACIÓN DIGITAL'
$ws.Range("E30").Value = 'Puertas, Edwin'
$ws.Range("F30").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57202285682'
$ws.Range("G30").Value = 'Puertas, Edwin (57202285682)'
$ws.Range("H30").Value = 4
$ws.Range("I30").Value = 4
$ws.Range("J30").Value = 1
$ws.Range("K30").Value = 1
$ws.Range("L30").Value = 1
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("R30").Value = 3
$ws.Range("S30").Value = 12
$ws.Range("U30").Value = 1
$ws.Range("V30").Value = 1
$ws.Range("W30").Value = 1
$ws.Range("X30").Value = 10
$ws.Range("AA30").Value = 0
$ws.Range("AC30").Value = 2
$ws.Range("AD30").Value = 3
$ws.Range("AE30").Value = 3
$ws.Range("AF30").Value = 3
$ws.Range("AG30").Value = 15
$ws.Range("AH30").Value = 3
$ws.Range("AI30").Value = 13
$ws.Range("AJ30").Value = 1

# Row 31
$ws.Range("A31").Value = "'57205400052"
$ws.Range("B31").Value = 'T00012602'
$ws.Range("C31").Value = 'JENIFER YORIS VASQUEZ AGUILAR'
$ws.Range("E31").Value = 'Vásquez, Jenifer'
$ws.Range("F31").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57205400052'
$ws.Range("G31").Value = 'Vásquez, Jenifer (57205400052)'
$ws.Range("L31").Value = 0
$ws.Range("R31").Value = 0
$ws.Range("S31").Value = 2
$ws.Range("X31").Value = 4
$ws.Range("AB31").Value = 0
$ws.Range("AC31").Value = 0
$ws.Range("AD31").Value = 0
$ws.Range("AE31").Value = 0
$ws.Range("AF31").Value = 0
$ws.Range("AG31").Value = 2
$ws.Range("AI31").Value = 4
$ws.Range("AJ31").Value = 0
$ws.Range("AK31").Value = 0

# Row 32
$ws.Range("A32").Value = "'57197327858"
$ws.Range("B32").Value = 'T00000049'
$ws.Range("C32").Value = 'OSCAR ACEVEDO PATIÑO'
$ws.Range("E32").Value = 'Acevedo-Patĩno, Oscar'
$ws.Range("F32").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57197327858'
$ws.Range("G32").Value = 'Acevedo, Oscar (57197327858)'
$ws.Range("L32").Value = 1
$ws.Range("R32").Value = 2
$ws.Range("S32").Value = 1
$ws.Range("AC32").Value = 1
$ws.Range("AD32").Value = 1
$ws.Range("AE32").Value = 1
$ws.Range("AF32").Value = 2
$ws.Range("AG32").Value = 3
$ws.Range("AK32").Value = 2

# Row 33
$ws.Range("A33").Value = "'56682785300"
$ws.Range("B33").Value = 'T00022128'
$ws.Range("C33").Value = 'LUZ ALEJANDRA MAGRE COLORADO'
$ws.Range("D33").Value = 'ESCUELA DE INGENIERÍA, ARQUITECTURA Y DISEÑO'
$ws.Range("E33").Value = 'Magre, Luz Alejandra'
$ws.Range("F33").Value = 'https://www.scopus.com/authid/detail.uri?authorId=56682785300'
$ws.Range("G33").Value = 'Magre Colorado, Luz A. (56682785300)'
$ws.Range("I33").Value = 1
$ws.Range("S33").Value = 0
$ws.Range("T33").Value = 0
$ws.Range("AB33").Value = 1
$ws.Range("AG33").Value = 0
$ws.Range("AH33").Value = 0
$ws.Range("AI33").Value = 0
$ws.Range("AJ33").Value = 1
$ws.Range("AK33").Value = 1

# Row 34
$ws.Range("A34").Value = "'57392556500"
$ws.Range("B34").Value = 'T00014731'
$ws.Range("C34").Value = 'LINA MARGARITA MARRUGO SALAS'
$ws.Range("D34").Value = 'ESCUELA DE NEGOCIOS, LEYES Y SOCIEDAD'
$ws.Range("E34").Value = 'Marrugo-Salas, Lina Margarita'
$ws.Range("F34").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57192271754'
$ws.Range("G34").Value = 'Marrugo-Salas, Lina (57392556500)'
$ws.Range("I34").Value = 2
$ws.Range("J34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("S34").Value = 4
$ws.Range("T34").Value = 1
$ws.Range("Y34").Value = 0
$ws.Range("AA34").Value = 0
$ws.Range("AD34").Value = 0
$ws.Range("AE34").Value = 0
$ws.Range("AF34").Value = 0
$ws.Range("AG34").Value = 4
$ws.Range("AH34").Value = 1
$ws.Range("AI34").Value = 1
$ws.Range("AJ34").Value = 0
$ws.Range("AK34").Value = 0

# Row 35
$ws.Range("A35").Value = "'57193252278"
$ws.Range("B35").Value = 'T00021661'
$ws.Range("C35").Value = 'JAIRO HUMBERTO CABRERA TOVAR'
$ws.Range("E35").Value = 'Cabrera, Jairo'
$ws.Range("F35").Value = 'https://scopus.utb.elogim.com/authid/detail.uri?authorId=57193252278'
$ws.Range("G35").Value = 'Cabrera, Jairo (57193252278)'
$ws.Range("S35").Value = 0
$ws.Range("X35").Value = 2
$ws.Range("AB35").Value = 1
$ws.Range("AG35").Value = 0
$ws.Range("AI35").Value = 2
$ws.Range("AJ35").Value = 1
$ws.Range("AK35").Value = 1

# Row 36
$ws.Range("A36").Value = "'57220927199"
$ws.Range("B36").Value = 'T00010915'
$ws.Range("C36").Value = 'JAIRO ENRIQUE SERRANO CASTAÑEDA'
$ws.Range("D36").Value = 'ESCUELA DE TRANSFORMACIÓN DIGITAL'
$ws.Range("E36").Value = 'Serrano, Jairo E.'
$ws.Range("F36").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57220927199'
$ws.Range("G36").Value = 'Serrano, Jairo E. (57220927199)'
$ws.Range("I36").Value = 2
$ws.Range("R36").Value = 1
$ws.Range("W36").Value = 2
$ws.Range("X36").Value = 0
$ws.Range("AF36").Value = 1
$ws.Range("AG36").Value = 1
$ws.Range("AH36").Value = 2

# Row 37
$ws.Range("A37").Value = "'57203321995"
$ws.Range("B37").Value = 'T00021700'
$ws.Range("C37").Value = 'HERNANDO RAFAEL ALTAMAR MERCADO'
$ws.Range("D37").Value = 'CIENCIAS BÁSICAS'
$ws.Range("E37").Value = 'Altamar-Mercado, Hernando'
$ws.Range("F37").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57203321995'
$ws.Range("G37").Value = 'Altamar-Mercado, Hernando (57203321995)'
$ws.Range("I37").Value = 0
$ws.Range("R37").Value = 0
$ws.Range("S37").Value = 3
$ws.Range("W37").Value = 1
$ws.Range("AB37").Value = 0
$ws.Range("AF37").Value = 0
$ws.Range("AG37").Value = 3
$ws.Range("AH37").Value = 1
$ws.Range("AI37").Value = 1
$ws.Range("AJ37").Value = 0
$ws.Range("AK37").Value = 0

# Row 38
$ws.Range("A38").Value = "'55258973100"
$ws.Range("B38").Value = 'T00059175'
$ws.Range("C38").Value = 'CLAUDIA PATRICIA DIAZ MENDOZA'
$ws.Range("D38").Value = 'ESCUELA DE INGENIERÍA, ARQUITECTURA Y DISEÑO'
$ws.Range("E38").Value = 'Díaz-Mendoza, Claudia'
$ws.Range("F38").Value = 'https://www.scopus.com/authid/detail.uri?authorId=55258973100'
$ws.Range("G38").Value = 'Díaz-Mendoza, Claudia (55258973100)'
$ws.Range("I38").Value = 2
$ws.Range("J38").Value = 2
$ws.Range("K38").Value = 1
$ws.Range("R38").Value = 1
$ws.Range("S38").Value = 1
$ws.Range("T38").Value = 1
$ws.Range("W38").Value = 0
$ws.Range("X38").Value = 1
$ws.Range("AB38").Value = 1
$ws.Range("AD38").Value = 3
$ws.Range("AE38").Value = 3
$ws.Range("AF38").Value = 1
$ws.Range("AG38").Value = 2
$ws.Range("AI38").Value = 2
$ws.Range("AJ38").Value = 1
$ws.Range("AK38").Value = 1

# Row 39
$ws.Range("A39").Value = "'57190688459"
$ws.Range("B39").Value = 'T00019348'
$ws.Range("C39").Value = 'ALBERTO PATIÑO VANEGAS'
$ws.Range("D39").Value = 'CIENCIAS BÁSICAS'
$ws.Range("E39").Value = 'Patiño, Alberto'
$ws.Range("F39").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57190688459'
$ws.Range("G39").Value = 'Patiño-Vanegas, Alberto (57190688459)'
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 0
$ws.Range("R39").Value = 0
$ws.Range("S39").Value = 5
$ws.Range("T39").Value = 0
$ws.Range("W39").Value = 1
$ws.Range("X39").Value = 4
$ws.Range("AB39").Value = 0
$ws.Range("AC39").Value = 1
$ws.Range("AD39").Value = 1
$ws.Range("AE39").Value = 1
$ws.Range("AF39").Value = 0
$ws.Range("AG39").Value = 5
$ws.Range("AI39").Value = 5
$ws.Range("AJ39").Value = 0

# Row 40
$ws.Range("A40").Value = "'57223851529"
$ws.Range("B40").Value = 'T00040584'
$ws.Range("C40").Value = 'ELSY MERCEDES DOMINGUEZ DE LA OSSA'
$ws.Range("D40").Value = 'ESCUELA DE NEGOCIOS, LEYES Y SOCIEDAD'
$ws.Range("E40").Value = 'de la Ossa, Elsy Domínguez'
$ws.Range("F40").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57223851529'
$ws.Range("G40").Value = 'Domínguez-De la Ossa, Elsy (57223851529)'
$ws.Range("I40").Value = 3
$ws.Range("K40").Value = 0
$ws.Range("P40").Value = 1
$ws.Range("Q40").Value = 1
$ws.Range("R40").Value = 2
$ws.Range("S40").Value = 0
$ws.Range("W40").Value = 0
$ws.Range("X40").Value = 0
$ws.Range("Y40").Value = 1
$ws.Range("AA40").Value = 2
$ws.Range("AC40").Value = 0
$ws.Range("AD40").Value = 0
$ws.Range("AE40").Value = 0
$ws.Range("AF40").Value = 4
$ws.Range("AG40").Value = 4
$ws.Range("AH40").Value = 0
$ws.Range("AI40").Value = 0
$ws.Range("AJ40").Value = 3
$ws.Range("AK40").Value = 3

# Row 41
$ws.Range("A41").Value = "'58068069000"
$ws.Range("B41").Value = 'T00000054'
$ws.Range("C41").Value = 'WILLIAM ARELLANO CARTAGENA'
$ws.Range("E41").Value = 'Arellano-Cartagena, William'
$ws.Range("F41").Value = 'https://www.scopus.com/authid/detail.uri?authorId=58068069000'
$ws.Range("G41").Value = 'Arellano-Cartagena, William (58068069000)'
$ws.Range("I41").Value = 2
$ws.Range("K41").Value = 0
$ws.Range("R41").Value = 1
$ws.Range("V41").Value = 1
$ws.Range("AD41").Value = 0
$ws.Range("AE41").Value = 0
$ws.Range("AF41").Value = 1
$ws.Range("AG41").Value = 1
$ws.Range("AH41").Value = 1
$ws.Range("AI41").Value = 1

# Row 42
$ws.Range("A42").Value = "'57218294431"
$ws.Range("B42").Value = 'T00051182'
$ws.Range("C42").Value = 'SERGIO VILLAR SALINAS'
$ws.Range("D42").Value = 'ESCUELA DE INGENIERÍA, ARQUITECTURA Y DISEÑO'
$ws.Range("E42").Value = 'Villar-Salinas, Sergio'
$ws.Range("F42").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57218294431'
$ws.Range("G42").Value = 'Villar-Salinas, Sergio (57218294431)'
$ws.Range("I42").Value = 2
$ws.Range("L42").Value = 0
$ws.Range("O42").Value = 1
$ws.Range("P42").Value = 1
$ws.Range("Q42").Value = 0
$ws.Range("AA42").Value = 0
$ws.Range("AD42").Value = 0
$ws.Range("AE42").Value = 0
$ws.Range("AF42").Value = 2
$ws.Range("AG42").Value = 2
$ws.Range("AJ42").Value = 0
$ws.Range("AK42").Value = 0

# Row 43
$ws.Range("A43").Value = "'55783129400"
$ws.Range("B43").Value = 'T00000064'
$ws.Range("C43").Value = 'MOISES RAMON QUINTANA ALVAREZ'
$ws.Range("E43").Value = 'Quintana, Moisés'
$ws.Range("F43").Value = 'https://www.scopus.com/authid/detail.uri?authorId=55783129400'
$ws.Range("G43").Value = 'Quintana, Moisés (55783129400)'
$ws.Range("I43").Value = 0
$ws.Range("S43").Value = 0
$ws.Range("T43").Value = 0
$ws.Range("X43").Value = 0
$ws.Range("Z43").Value = 0
$ws.Range("AB43").Value = 0
$ws.Range("AG43").Value = 0
$ws.Range("AH43").Value = 0
$ws.Range("AI43").Value = 0
$ws.Range("AJ43").Value = 0
$ws.Range("AK43").Value = 0

# Row 44
$ws.Range("A44").Value = "'57918628600"
$ws.Range("B44").Value = 'T00057400'
$ws.Range("C44").Value = 'JORGE LUIS VILLALBA ACEVEDO'
$ws.Range("D44").Value = 'CIENCIAS BÁSICAS'
$ws.Range("E44").Value = 'Villalba-Acevedo, Jorge Luis'
$ws.Range("F44").Value = 'https://scopus.utb.elogim.com/authid/detail.uri?authorId=57918628600'
$ws.Range("G44").Value = 'Villalba-Acevedo, Jorge Luis (57918628600)'
$ws.Range("R44").Value = 0
$ws.Range("S44").Value = 2
$ws.Range("T44").Value = 1
$ws.Range("V44").Value = 0
$ws.Range("X44").Value = 1
$ws.Range("AB44").Value = 1
$ws.Range("AF44").Value = 0
$ws.Range("AG44").Value = 2
$ws.Range("AI44").Value = 2
$ws.Range("AJ44").Value = 2
$ws.Range("AK44").Value = 2

# Row 45
$ws.Range("A45").Value = "'57221229836"
$ws.Range("B45").Value = 'T00015391'
$ws.Range("C45").Value = 'MARIA FERNANDA MEDINA REYES'
$ws.Range("D45").Value = 'ESCUELA DE TRANSFORMACIÓN DIGITAL'
$ws.Range("E45").Value = 'Medina, M. F.'
$ws.Range("F45").Value = 'https://scopus.utb.elogim.com/authid/detail.uri?authorId=57221229836'
$ws.Range("G45").Value = 'Medina-Reyes, María Fernanda (57221229836)'
$ws.Range("I45").Value = 1
$ws.Range("L45").Value = 1
$ws.Range("Q45").Value = 1
$ws.Range("AA45").Value = 1
$ws.Range("AD45").Value = 1
$ws.Range("AE45").Value = 1
$ws.Range("AF45").Value = 1
$ws.Range("AG45").Value = 1
$ws.Range("AJ45").Value = 1
$ws.Range("AK45").Value = 1

# Row 46
$ws.Range("A46").Value = "'57204842254"
$ws.Range("B46").Value = 'T00019354'
$ws.Range("C46").Value = 'TANIA ISABEL JIMENEZ CASTILLA'
$ws.Range("D46").Value = 'ESCUELA DE NEGOCIOS, LEYES Y SOCIEDAD'
$ws.Range("E46").Value = 'Jiménez-Castilla, Tania'
$ws.Range("F46").Value = 'https://www.scopus.com/authid/detail.uri?authorId=57204842254'
$ws.Range("G46").Value = 'Jiménez-Castilla, Tania (57204842254)'
$ws.Range("I46").Value = 1
$ws.Range("K46").Value = 1
$ws.Range("O46").Value = 0
$ws.Range("P46").Value = 0
$ws.Range("Z46").Value = 1
$ws.Range("AD46").Value = 1
$ws.Range("AE46").Value = 1
$ws.Range("AF46").Value = 0
$ws.Range("AG46").Value = 0
$ws.Range("AJ46").Value = 1
$ws.Range("AK46").Value = 1
